$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 4.616372666666667
$ws.Range("H2").Value = 13.849118
$ws.Range("I2").Value = 0.0183283362562958
$ws.Range("J2").Value = 0.01832833625629581
$ws.Range("M2").Value = 0.1404576666666667
$ws.Range("N2").Value = 0.421373
$ws.Range("O2").Value = 0.183302244487095
$ws.Range("P2").Value = 0.1833022444870949
$ws.Range("Q2").Value = 0.6484049332237779
$ws.Range("R2").Value = 5.835644399014
$ws.Range("S2").Value = 0.00335962517349322
$ws.Range("T2").Value = 0.00335962517349322

# Row 3
$ws.Range("G3").Value = 4.616372666666667
$ws.Range("H3").Value = 13.849118
$ws.Range("I3").Value = 0.0183283362562958
$ws.Range("J3").Value = 0.01832833625629581
$ws.Range("M3").Value = 0.5312866666666666
$ws.Range("O3").Value = 0.6933479729318232
$ws.Range("P3").Value = 0.6933479729318232
$ws.Range("Q3").Value = 2.452617246164444
$ws.Range("R3").Value = 22.07355521548
$ws.Range("S3").Value = 0.01270791479051554
$ws.Range("T3").Value = 0.01270791479051554

# Row 4
$ws.Range("G4").Value = 4.616372666666667
$ws.Range("H4").Value = 13.849118
$ws.Range("I4").Value = 0.0183283362562958
$ws.Range("J4").Value = 0.01832833625629581
$ws.Range("O4").Value = 0.1233497825810819
$ws.Range("P4").Value = 0.1233497825810819
$ws.Range("Q4").Value = 0.4363318504988889
$ws.Range("R4").Value = 3.92698665449
$ws.Range("S4").Value = 0.002260796292287047
$ws.Range("T4").Value = 0.002260796292287048

# Row 5
$ws.Range("I5").Value = 0.943783113604627
$ws.Range("J5").Value = 0.9437831136046271
$ws.Range("M5").Value = 0.1404576666666667
$ws.Range("N5").Value = 0.421373
$ws.Range("O5").Value = 0.183302244487095
$ws.Range("P5").Value = 0.1833022444870949
$ws.Range("Q5").Value = 33.38838933317422
$ws.Range("R5").Value = 300.495503998568
$ws.Range("S5").Value = 0.172997563032747
$ws.Range("T5").Value = 0.172997563032747

# Row 6
$ws.Range("I6").Value = 0.943783113604627
$ws.Range("J6").Value = 0.9437831136046271
$ws.Range("M6").Value = 0.5312866666666666
$ws.Range("O6").Value = 0.6933479729318232
$ws.Range("P6").Value = 0.6933479729318232
$ws.Range("Q6").Value = 126.2929001681955
$ws.Range("S6").Value = 0.6543701087050526
$ws.Range("T6").Value = 0.6543701087050527

# Row 7
$ws.Range("I7").Value = 0.943783113604627
$ws.Range("J7").Value = 0.9437831136046271
$ws.Range("O7").Value = 0.1233497825810819
$ws.Range("P7").Value = 0.1233497825810819
$ws.Range("S7").Value = 0.1164154418668272
$ws.Range("T7").Value = 0.1164154418668272

# Row 8
$ws.Range("I8").Value = 0.03788855013907712
$ws.Range("J8").Value = 0.03788855013907712
$ws.Range("M8").Value = 0.1404576666666667
$ws.Range("N8").Value = 0.421373
$ws.Range("O8").Value = 0.183302244487095
$ws.Range("P8").Value = 0.1833022444870949
$ws.Range("Q8").Value = 1.340390228514889
$ws.Range("R8").Value = 12.063512056634
$ws.Range("S8").Value = 0.006945056280854669
$ws.Range("T8").Value = 0.006945056280854669

# Row 9
$ws.Range("I9").Value = 0.03788855013907712
$ws.Range("J9").Value = 0.03788855013907712
$ws.Range("M9").Value = 0.5312866666666666
$ws.Range("O9").Value = 0.6933479729318232
$ws.Range("P9").Value = 0.6933479729318232
$ws.Range("Q9").Value = 5.070078931542222
$ws.Range("R9").Value = 45.63071038387999
$ws.Range("S9").Value = 0.02626994943625487
$ws.Range("T9").Value = 0.02626994943625487

# Row 10
$ws.Range("I10").Value = 0.03788855013907712
$ws.Range("J10").Value = 0.03788855013907712
$ws.Range("O10").Value = 0.1233497825810819
$ws.Range("P10").Value = 0.1233497825810819
$ws.Range("S10").Value = 0.004673544421967581
$ws.Range("T10").Value = 0.004673544421967582
